# Updates the cryptocurrency price/volume snapshot (columns D and E) for
# rows 2-51 on the active sheet, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '60.801.30'; ForceText = $false },
    @{ Cell = "E2"; Value = '  -1.60%  '; ForceText = $false },
    @{ Cell = "D3"; Value = '2.907.61'; ForceText = $false },
    @{ Cell = "E3"; Value = '  -2.55%  '; ForceText = $false },
    @{ Cell = "E4"; Value = '  +0.04%  '; ForceText = $false },
    @{ Cell = "D5"; Value = '527.79'; ForceText = $true },
    @{ Cell = "E5"; Value = '  -2.45%  '; ForceText = $false },
    @{ Cell = "D6"; Value = '143.84'; ForceText = $true },
    @{ Cell = "E6"; Value = '  -5.44%  '; ForceText = $false },
    @{ Cell = "D7"; Value = '0.999'; ForceText = $true },
    @{ Cell = "E7"; Value = '  +0.03%  '; ForceText = $false },
    @{ Cell = "D8"; Value = '0.546'; ForceText = $true },
    @{ Cell = "E8"; Value = '  -4.10%  '; ForceText = $false },
    @{ Cell = "D9"; Value = '2.915.60'; ForceText = $false },
    @{ Cell = "E9"; Value = '  -2.60%  '; ForceText = $false },
    @{ Cell = "E10"; Value = '  -4.54%  '; ForceText = $false },
    @{ Cell = "E11"; Value = '  -1.64%  '; ForceText = $false },
    @{ Cell = "D12"; Value = '0.358'; ForceText = $true },
    @{ Cell = "E12"; Value = '  -2.91%  '; ForceText = $false },
    @{ Cell = "D13"; Value = '3.413.50'; ForceText = $false },
    @{ Cell = "E13"; Value = '  -2.55%  '; ForceText = $false },
    @{ Cell = "E14"; Value = '  +2.97%  '; ForceText = $false },
    @{ Cell = "D15"; Value = '60.747.28'; ForceText = $false },
    @{ Cell = "E15"; Value = '  -1.72%  '; ForceText = $false },
    @{ Cell = "D16"; Value = '22.53'; ForceText = $true },
    @{ Cell = "E16"; Value = '  -5.89%  '; ForceText = $false },
    @{ Cell = "D17"; Value = '2.896.48'; ForceText = $false },
    @{ Cell = "E17"; Value = '  -2.99%  '; ForceText = $false },
    @{ Cell = "E18"; Value = '  -4.19%  '; ForceText = $false },
    @{ Cell = "D19"; Value = '4.96'; ForceText = $true },
    @{ Cell = "E19"; Value = '  -3.94%  '; ForceText = $false },
    @{ Cell = "D20"; Value = '11.58'; ForceText = $true },
    @{ Cell = "E20"; Value = '  -3.90%  '; ForceText = $false },
    @{ Cell = "D21"; Value = '353.21'; ForceText = $true },
    @{ Cell = "E21"; Value = '  -7.32%  '; ForceText = $false },
    @{ Cell = "D22"; Value = '6.51'; ForceText = $true },
    @{ Cell = "E22"; Value = '  -2.98%  '; ForceText = $false },
    @{ Cell = "E23"; Value = '  +0.05%  '; ForceText = $false },
    @{ Cell = "E24"; Value = '  +1.46%  '; ForceText = $false },
    @{ Cell = "D25"; Value = '65.00'; ForceText = $true },
    @{ Cell = "E25"; Value = '  -1.53%  '; ForceText = $false },
    @{ Cell = "E26"; Value = '  -4.19%  '; ForceText = $false },
    @{ Cell = "E27"; Value = '  -6.60%  '; ForceText = $false },
    @{ Cell = "E28"; Value = '  +0.19%  '; ForceText = $false },
    @{ Cell = "D29"; Value = '7.85'; ForceText = $true },
    @{ Cell = "E29"; Value = '  -3.51%  '; ForceText = $false },
    @{ Cell = "D30"; Value = '0.0₃0855'; ForceText = $false },
    @{ Cell = "E30"; Value = '  -9.25%  '; ForceText = $false },
    @{ Cell = "E31"; Value = '  +0.01%  '; ForceText = $false },
    @{ Cell = "D32"; Value = '1.68'; ForceText = $true },
    @{ Cell = "E32"; Value = '  -1.76%  '; ForceText = $false },
    @{ Cell = "D33"; Value = '19.59'; ForceText = $true },
    @{ Cell = "E33"; Value = '  -4.36%  '; ForceText = $false },
    @{ Cell = "D34"; Value = '153.25'; ForceText = $true },
    @{ Cell = "E34"; Value = '  -4.03%  '; ForceText = $false },
    @{ Cell = "D35"; Value = '4.38'; ForceText = $true },
    @{ Cell = "E35"; Value = '  -4.14%  '; ForceText = $false },
    @{ Cell = "D36"; Value = '5.57'; ForceText = $true },
    @{ Cell = "E36"; Value = '  -5.89%  '; ForceText = $false },
    @{ Cell = "D37"; Value = '0.995'; ForceText = $true },
    @{ Cell = "E37"; Value = '  -7.10%  '; ForceText = $false },
    @{ Cell = "E38"; Value = '  -5.77%  '; ForceText = $false },
    @{ Cell = "D39"; Value = '37.51'; ForceText = $true },
    @{ Cell = "E39"; Value = '  -0.19%  '; ForceText = $false },
    @{ Cell = "E40"; Value = '  -4.73%  '; ForceText = $false },
    @{ Cell = "D41"; Value = '3.72'; ForceText = $true },
    @{ Cell = "E41"; Value = '  -4.83%  '; ForceText = $false },
    @{ Cell = "D42"; Value = '2.291.67'; ForceText = $false },
    @{ Cell = "E42"; Value = '  -5.28%  '; ForceText = $false },
    @{ Cell = "D43"; Value = '0.651'; ForceText = $true },
    @{ Cell = "E43"; Value = '  -3.23%  '; ForceText = $false },
    @{ Cell = "D44"; Value = '0.0582'; ForceText = $true },
    @{ Cell = "E44"; Value = '  -1.28%  '; ForceText = $false },
    @{ Cell = "D45"; Value = '20.35'; ForceText = $true },
    @{ Cell = "E45"; Value = '  -7.46%  '; ForceText = $false },
    @{ Cell = "E46"; Value = '  +0.06%  '; ForceText = $false },
    @{ Cell = "D47"; Value = '4.96'; ForceText = $true },
    @{ Cell = "E47"; Value = '  -3.25%  '; ForceText = $false },
    @{ Cell = "D48"; Value = '0.0238'; ForceText = $true },
    @{ Cell = "E48"; Value = '  -2.69%  '; ForceText = $false },
    @{ Cell = "D49"; Value = '10.33'; ForceText = $true },
    @{ Cell = "E49"; Value = '  -0.91%  '; ForceText = $false },
    @{ Cell = "D50"; Value = '0.0917'; ForceText = $true },
    @{ Cell = "E50"; Value = '  -3.69%  '; ForceText = $false },
    @{ Cell = "D51"; Value = '18.41'; ForceText = $true },
    @{ Cell = "E51"; Value = '  -7.08%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Values such as "527.79" would otherwise be auto-converted to a
        # number by Excel; force text, assign, then drop the now-redundant
        # explicit "@" format so the cell style matches a plain text cell.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
